$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verifications")

# Add new verification entry text for Costs Decision Cover Letter in A17
$ws.Range("A17").Value = "Costs Decision Cover Letter"

# Update the active selection to A14
$ws.Range("A14").Select()
